$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 6: B6 "+" -> "-", C6 35 -> 36 ---
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = 36

# --- Append new rows 60:70 (word "g" = 1+1+1+1+1) ---

# Copy column-A formatting (bold, centered, bordered) from the last
# existing data row down across the new rows before writing values.
$ws.Range("A59").Copy()
$ws.Range("A60:A70").PasteSpecial(-4122)

$newRows = @(
    @(60, 58, "g", 81),
    @(61, 59, "=", 46),
    @(62, 60, "1", 80),
    @(63, 61, "+", 35),
    @(64, 62, "1", 80),
    @(65, 63, "+", 35),
    @(66, 64, "1", 80),
    @(67, 65, "+", 35),
    @(68, 66, "1", 80),
    @(69, 67, "+", 35),
    @(70, 68, "1", 80)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $c = $row[3]
    $ws.Range("A$r").Value = $a
    # Column B holds tokenizer symbols as literal TEXT (even things that
    # look like numbers, e.g. "1", or that look like formulas, e.g. "=").
    # Prefix with an apostrophe so Excel stores them as text instead of
    # auto-converting to a number / parsing as a formula.
    $ws.Range("B$r").Value = "'" + $b
    $ws.Range("C$r").Value = $c
}

# The apostrophe text-prefix marks the cells with a "quote prefix" style;
# reset back to the default Normal style so formatting matches the other
# plain-text cells in the column (no visible leading apostrophe either way).
$ws.Range("B60:B70").Style = "Normal"
